$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '302.75'
$ws.Range('D2').NumberFormat = "General"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '1.91%'
$ws.Range('E2').NumberFormat = "General"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '31.92'
$ws.Range('D3').NumberFormat = "General"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '2.03%'
$ws.Range('E3').NumberFormat = "General"

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.115'
$ws.Range('D4').NumberFormat = "General"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '0.50%'
$ws.Range('E4').NumberFormat = "General"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07820'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-2.52%'
$ws.Range('E5').NumberFormat = "General"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '2.105'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '-22.00%'
$ws.Range('E6').NumberFormat = "General"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '7.813'
$ws.Range('D7').NumberFormat = "General"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '0.30%'
$ws.Range('E7').NumberFormat = "General"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.794'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-0.36%'
$ws.Range('E8').NumberFormat = "General"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9218'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-0.65%'
$ws.Range('E9').NumberFormat = "General"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1757'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '0.23%'
$ws.Range('E10').NumberFormat = "General"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07693'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '5.26%'
$ws.Range('E11').NumberFormat = "General"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08904'
$ws.Range('D12').NumberFormat = "General"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-1.31%'
$ws.Range('E12').NumberFormat = "General"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03142'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '4.15%'
$ws.Range('E13').NumberFormat = "General"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-0.11%'
$ws.Range('E14').NumberFormat = "General"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001517'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '1.63%'
$ws.Range('E15').NumberFormat = "General"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.005870'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-1.06%'
$ws.Range('E16').NumberFormat = "General"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.450'
$ws.Range('D17').NumberFormat = "General"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-2.38%'
$ws.Range('E17').NumberFormat = "General"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.261'
$ws.Range('D18').NumberFormat = "General"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '0.69%'
$ws.Range('E18').NumberFormat = "General"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3293'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.91%'
$ws.Range('E19').NumberFormat = "General"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.1331'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '-0.20%'
$ws.Range('E20').NumberFormat = "General"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.173'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '0.85%'
$ws.Range('E21').NumberFormat = "General"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.1794'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '8.95%'
$ws.Range('E22').NumberFormat = "General"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04584'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '-0.30%'
$ws.Range('E23').NumberFormat = "General"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.001238'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-0.29%'
$ws.Range('E24').NumberFormat = "General"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '1.04%'
$ws.Range('E25').NumberFormat = "General"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001250'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '4.25%'
$ws.Range('E26').NumberFormat = "General"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01771'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '0.35%'
$ws.Range('E39').NumberFormat = "General"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04771'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '6.07%'
$ws.Range('E40').NumberFormat = "General"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007008'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '2.68%'
$ws.Range('E41').NumberFormat = "General"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1368'
$ws.Range('D42').NumberFormat = "General"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '1.49%'
$ws.Range('E42').NumberFormat = "General"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.002190'
$ws.Range('D43').NumberFormat = "General"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '2.61%'
$ws.Range('E43').NumberFormat = "General"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.01073'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '8.99%'
$ws.Range('E44').NumberFormat = "General"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00006268'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-3.65%'
$ws.Range('E45').NumberFormat = "General"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '0.09%'
$ws.Range('E46').NumberFormat = "General"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.003561'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-59.28%'
$ws.Range('E47').NumberFormat = "General"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.8073'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-1.62%'
$ws.Range('E48').NumberFormat = "General"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '0.09%'
$ws.Range('E49').NumberFormat = "General"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0002000'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.09%'
$ws.Range('E50').NumberFormat = "General"

